$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "28.470.59"
Set-TextValue $ws.Range("E2") "  -5.37%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.844.46"
Set-TextValue $ws.Range("E3") "  -5.39%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.005"
Set-TextValue $ws.Range("E4") "  -0.39%  "

# Row 5
Set-TextValue $ws.Range("D5") "334.66"
Set-TextValue $ws.Range("E5") "  +1.90%  "

# Row 6
Set-TextValue $ws.Range("E6") "  -0.45%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.4628"
Set-TextValue $ws.Range("E7") "  -4.84%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.3885"
Set-TextValue $ws.Range("E8") "  -5.25%  "

# Row 9
Set-TextValue $ws.Range("D9") "46.00"
Set-TextValue $ws.Range("E9") "  -3.65%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.07876"
Set-TextValue $ws.Range("E10") "  -4.53%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.9719"
Set-TextValue $ws.Range("E11") "  -4.80%  "

# Row 12
Set-TextValue $ws.Range("D12") "22.05"
Set-TextValue $ws.Range("E12") "  -8.44%  "

# Row 13
Set-TextValue $ws.Range("D13") "1.874.72"
Set-TextValue $ws.Range("E13") "  -4.25%  "

# Row 14
Set-TextValue $ws.Range("D14") "5.800"
Set-TextValue $ws.Range("E14") "  -5.41%  "

# Row 15
Set-TextValue $ws.Range("D15") "6.944"
Set-TextValue $ws.Range("E15") "  -5.69%  "

# Row 16
Set-TextValue $ws.Range("D16") "0.06896"
Set-TextValue $ws.Range("E16") "  +0.19%  "

# Row 17
Set-TextValue $ws.Range("E17") "  -0.40%  "

# Row 18
Set-TextValue $ws.Range("D18") "87.50"
Set-TextValue $ws.Range("E18") "  -4.82%  "

# Row 19
Set-TextValue $ws.Range("D19") "0.00001002"
Set-TextValue $ws.Range("E19") "  -3.80%  "

# Row 20
Set-TextValue $ws.Range("D20") "17.00"
Set-TextValue $ws.Range("E20") "  -4.92%  "

# Row 21
Set-TextValue $ws.Range("D21") "1.003"
Set-TextValue $ws.Range("E21") "  -0.77%  "

# Row 22
Set-TextValue $ws.Range("D22") "28.479.54"
Set-TextValue $ws.Range("E22") "  -4.99%  "

# Row 23
Set-TextValue $ws.Range("D23") "5.364"
Set-TextValue $ws.Range("E23") "  -6.12%  "

# Row 24
Set-TextValue $ws.Range("E24") "  -7.35%  "

# Row 25
Set-TextValue $ws.Range("D25") "2.171"
Set-TextValue $ws.Range("E25") "  -1.39%  "

# Row 26
Set-TextValue $ws.Range("D26") "2.081.63"
Set-TextValue $ws.Range("E26") "  -4.78%  "

# Row 27
Set-TextValue $ws.Range("D27") "153.53"
Set-TextValue $ws.Range("E27") "  -2.13%  "

# Row 28
Set-TextValue $ws.Range("D28") "19.35"
Set-TextValue $ws.Range("E28") "  -3.93%  "

# Row 29
Set-TextValue $ws.Range("D29") "5.908"
Set-TextValue $ws.Range("E29") "  -10.41%  "

# Row 30
Set-TextValue $ws.Range("D30") "1.990"
Set-TextValue $ws.Range("E30") "  -6.26%  "

# Row 31
Set-TextValue $ws.Range("D31") "117.29"
Set-TextValue $ws.Range("E31") "  -3.50%  "

# Row 32
Set-TextValue $ws.Range("D32") "0.9538"
Set-TextValue $ws.Range("E32") "  -7.10%  "

# Row 33
Set-TextValue $ws.Range("D33") "0.09361"
Set-TextValue $ws.Range("E33") "  -3.06%  "

# Row 34
Set-TextValue $ws.Range("D34") "5.337"
Set-TextValue $ws.Range("E34") "  -5.54%  "

# Row 35
Set-TextValue $ws.Range("E35") "  -2.83%  "

# Row 36
Set-TextValue $ws.Range("D36") "1.330"
Set-TextValue $ws.Range("E36") "  -6.70%  "

# Row 37
Set-TextValue $ws.Range("D37") "0.06054"
Set-TextValue $ws.Range("E37") "  -8.07%  "

# Row 38
Set-TextValue $ws.Range("D38") "0.02184"
Set-TextValue $ws.Range("E38") "  -5.21%  "

# Row 39
Set-TextValue $ws.Range("D39") "1.160"
Set-TextValue $ws.Range("E39") "  -4.54%  "

# Row 40
Set-TextValue $ws.Range("B40") "Frax"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue $ws.Range("D40") "1.001"
Set-TextValue $ws.Range("E40") "  -0.76%  "

# Row 41
Set-TextValue $ws.Range("B41") "TheSandbox"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D41") "0.5669"
Set-TextValue $ws.Range("E41") "  -5.43%  "

# Row 42
Set-TextValue $ws.Range("B42") "FraxShare"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D42") "7.614"
Set-TextValue $ws.Range("E42") "  -4.68%  "

# Row 43
Set-TextValue $ws.Range("B43") "Aptos"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D43") "10.04"
Set-TextValue $ws.Range("E43") "  -7.11%  "

# Row 44
Set-TextValue $ws.Range("B44") "Algorand"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D44") "0.1791"
Set-TextValue $ws.Range("E44") "  -3.71%  "

# Row 45
Set-TextValue $ws.Range("B45") "RenderToken"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D45") "2.387"
Set-TextValue $ws.Range("E45") "  -7.10%  "

# Row 46
Set-TextValue $ws.Range("B46") "WEMIXToken"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D46") "1.249"
Set-TextValue $ws.Range("E46") "  -0.06%  "

# Row 47
Set-TextValue $ws.Range("B47") "EnergySwap"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D47") "11.74"
Set-TextValue $ws.Range("E47") "  -6.04%  "

# Row 48
Set-TextValue $ws.Range("B48") "Decentraland"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D48") "0.5342"
Set-TextValue $ws.Range("E48") "  -4.55%  "

# Row 49
Set-TextValue $ws.Range("B49") "Cronos"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D49") "0.07056"
Set-TextValue $ws.Range("E49") "  -6.77%  "

# Row 50
Set-TextValue $ws.Range("B50") "NEARProtocol"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D50") "1.865"
Set-TextValue $ws.Range("E50") "  -6.33%  "

# Row 51
Set-TextValue $ws.Range("B51") "Quant"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D51") "112.86"
Set-TextValue $ws.Range("E51") "  -4.77%  "
